# Update Name of Algo
# Refresh the imputed "A" column values (delta-13C results) for the
# RandomForest run on the terrestrial_mammals / combination_2_ABCDE /
# A / 20 / seed3 dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = -22.39230000000002
    10  = -21.9772
    12  = -21.577
    18  = -22.31140000000002
    25  = -21.70719999999998
    37  = -19.51849999999999
    55  = -22.41980000000001
    68  = -21.53619999999999
    77  = -20.71259999999999
    78  = -20.57919999999997
    79  = -20.76549999999999
    80  = -19.4792
    81  = -21.82490000000001
    82  = -21.95519999999999
    84  = -22.07660000000001
    101 = -20.97429999999997
    102 = -19.61579999999999
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
